$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data import: participant names for the Exp1/Exp2 columns (E and F), rows 2-14 ---
$names = @(
    "Sumera",
    "Megan",
    "Tate",
    "Christine",
    "Mark",
    "Natalie",
    "Shruti",
    "Cooro",
    "Benita",
    "Agusta",
    "Ayush",
    "Sue",
    "Madeleine"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $names[$i]
    $ws.Cells.Item($row, 6).Value = $names[$i]
}

# --- Added analysis: new subject row 14 (subj 13), following the Exp1-first / no-disparity pattern ---
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Exp1"
$ws.Cells.Item(14, 3).Value = "Exp2"
$ws.Cells.Item(14, 4).Value = "no"

# --- Column widths to fit the newly added data ---
$ws.Columns.Item(4).ColumnWidth = 16.28515625
$ws.Range("E:F").ColumnWidth = 9.5703125

# --- Restore cursor/selection position ---
$ws.Range("R11").Select() | Out-Null
